$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells that are being rewritten keep their original text
# representation (Excel would otherwise auto-convert numeric-looking strings to numbers).
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D9").NumberFormat = "@"
$ws.Range("D11:D17").NumberFormat = "@"
$ws.Range("D19:D25").NumberFormat = "@"
$ws.Range("D27:D29").NumberFormat = "@"
$ws.Range("D31:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.297.38"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "3.696.05"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "683.39"
$ws.Range("E5").Value = "  -2.95%  "
$ws.Range("D6").Value = "162.91"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("D7").Value = "3.695.99"
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -5.29%  "
$ws.Range("E10").Value = "  -7.06%  "
$ws.Range("D11").Value = "7.23"
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  -5.93%  "
$ws.Range("D14").Value = "33.60"
$ws.Range("E14").Value = "  -7.04%  "
$ws.Range("D15").Value = "4.314.81"
$ws.Range("E15").Value = "  -3.09%  "
$ws.Range("D16").Value = "3.692.91"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "69.384.33"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "16.39"
$ws.Range("E19").Value = "  -4.78%  "
$ws.Range("D20").Value = "6.61"
$ws.Range("E20").Value = "  -7.55%  "
$ws.Range("D21").Value = "484.08"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").Value = "9.84"
$ws.Range("E22").Value = "  -6.63%  "
$ws.Range("D23").Value = "0.664"
$ws.Range("E23").Value = "  -8.52%  "
$ws.Range("D24").Value = "79.46"
$ws.Range("E24").Value = "  -7.00%  "
$ws.Range("D25").Value = "3.833.93"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("E26").Value = "  -9.63%  "
$ws.Range("D27").Value = "11.62"
$ws.Range("E27").Value = "  -3.49%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "9.51"
$ws.Range("E29").Value = "  -9.00%  "
$ws.Range("E30").Value = "  -11.89%  "
$ws.Range("D31").Value = "2.75"
$ws.Range("E31").Value = "  -10.47%  "
$ws.Range("D32").Value = "2.12"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").Value = "6.74"
$ws.Range("E33").Value = "  -8.25%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "26.81"
$ws.Range("E35").Value = "  -7.88%  "
$ws.Range("D36").Value = "0.163"
$ws.Range("E36").Value = "  -6.01%  "
$ws.Range("D37").Value = "3.659.32"
$ws.Range("E37").Value = "  -3.16%  "
$ws.Range("D38").Value = "8.56"
$ws.Range("E38").Value = "  -5.24%  "
$ws.Range("D39").Value = "6.08"
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("D40").Value = "0.0935"
$ws.Range("E40").Value = "  -7.53%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "0.958"
$ws.Range("E44").Value = "  -8.12%  "
$ws.Range("D45").Value = "161.12"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").Value = "48.29"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "2.85"
$ws.Range("E47").Value = "  -12.87%  "
$ws.Range("D48").Value = "393.48"
$ws.Range("E48").Value = "  -6.80%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "0.000277"
$ws.Range("E49").Value = "  -9.75%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.30"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "28.14"
$ws.Range("E51").Value = "  +1.82%  "
